$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 224, pushing existing rows 224..274 down to 225..275
$ws.Rows.Item(224).Insert()

# Populate the new row 224 with the new record
$ws.Range("A224").Value = 9
$ws.Range("B224").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C224").Value = "Metropolitana"
$ws.Range("D224").Value = 44951
$ws.Range("E224").Value = 13
$ws.Range("F224").Value = "Fruta"
$ws.Range("G224").Value = 100101
$ws.Range("H224").Value = "Berries"
$ws.Range("I224").Value = 100101001
$ws.Range("J224").Value = "Arándano (blue)"
$ws.Range("K224").Value = "Sin especificar"
$ws.Range("L224").Value = "Primera"
$ws.Range("M224").Value = 280
$ws.Range("N224").Value = 3000
$ws.Range("O224").Value = 3000
$ws.Range("P224").Value = 3000
$ws.Range("Q224").Value = "`$/bandeja 2 kilos"
$ws.Range("R224").Value = "Región de O'Higgins"
$ws.Range("S224").Value = 1500
$ws.Range("T224").Value = 2
